$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').NumberFormat = '@'
$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('B2').Style = 'Normal'
$ws.Range('C2').NumberFormat = '@'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('C2').Style = 'Normal'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.707.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('B3').NumberFormat = '@'
$ws.Range('B3').Value = 'Ethereum'
$ws.Range('B3').Style = 'Normal'
$ws.Range('C3').NumberFormat = '@'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('C3').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.305.50'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('B4').NumberFormat = '@'
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('B4').Style = 'Normal'
$ws.Range('C4').NumberFormat = '@'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('C4').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('B5').NumberFormat = '@'
$ws.Range('B5').Value = 'BNB'
$ws.Range('B5').Style = 'Normal'
$ws.Range('C5').NumberFormat = '@'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('C5').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = 'Solana'
$ws.Range('B6').Style = 'Normal'
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('C6').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'USDC'
$ws.Range('B7').Style = 'Normal'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('C7').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('B8').Style = 'Normal'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('C8').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.295.54'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'XRP'
$ws.Range('B9').Style = 'Normal'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('C9').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.572'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.62%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('B10').Style = 'Normal'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('C10').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.173'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.49%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'Cardano'
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.574'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.76%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.32%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000260'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.831.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.38'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.77%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '568.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -11.00%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.658.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'TRON'
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.307.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('B21').Style = 'Normal'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('C21').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.15%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Polygon'
$ws.Range('B22').Style = 'Normal'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('C22').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.885'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('B23').Style = 'Normal'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C23').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.66'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.28%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -8.73%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'LEO'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.82%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.44%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.59'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.44%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -8.18%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '555.38'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.14%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Cosmos'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.75'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Maker'
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.759.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.05%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Hedera'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.102'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.89%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Dai'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'OKB'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '55.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.75%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.28'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.125'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Stacks'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.11'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.06%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'PEPE'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₃0671'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -8.06%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.54'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -6.09%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.329'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'VeChain'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0404'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.20%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'CoreDAO'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -10.87%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Stellar'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.125'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.44%  '
$ws.Range('E51').Style = 'Normal'

Write-Output "done"